$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs / Nrg4 / Erbb4 / MuSCs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.352346
$ws.Range("H2").Value = 2.704692
$ws.Range("I2").Value = 0.1602206213737441
$ws.Range("J2").Value = 0.1162353457488538
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0108025
$ws.Range("N2").Value = 0.021605
$ws.Range("Q2").Value = 0.014608717665
$ws.Range("R2").Value = 0.05843487066
$ws.Range("S2").Value = 0.1602206213737441
$ws.Range("T2").Value = 0.1162353457488538

# Row 3 (FAPs / Nrg4 / Erbb4 / MuSCs)
$ws.Range("I3").Value = 0.01228059616519859
$ws.Range("J3").Value = 0.01336381667688838
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.0108025
$ws.Range("N3").Value = 0.021605
$ws.Range("Q3").Value = 0.001119729536666667
$ws.Range("R3").Value = 0.00671837722
$ws.Range("S3").Value = 0.01228059616519859
$ws.Range("T3").Value = 0.01336381667688838

# Row 4 (Inflammatory-Mac / Nrg4 / Erbb4 / MuSCs)
$ws.Range("G4").Value = 1.818045333333333
$ws.Range("H4").Value = 5.454136
$ws.Range("I4").Value = 0.2153948419948019
$ws.Range("J4").Value = 0.2343939286696121
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.0108025
$ws.Range("N4").Value = 0.021605
$ws.Range("Q4").Value = 0.01963943471333333
$ws.Range("R4").Value = 0.11783660828
$ws.Range("S4").Value = 0.2153948419948019
$ws.Range("T4").Value = 0.2343939286696121

# Row 5 (MuSCs / Nrg4 / Erbb4 / MuSCs)
$ws.Range("G5").Value = 0.7001250000000001
$ws.Range("H5").Value = 1.40025
$ws.Range("I5").Value = 0.08294804919694561
$ws.Range("J5").Value = 0.0601763686530047
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.0108025
$ws.Range("N5").Value = 0.021605
$ws.Range("Q5").Value = 0.007563100312500001
$ws.Range("R5").Value = 0.03025240125000001
$ws.Range("S5").Value = 0.08294804919694561
$ws.Range("T5").Value = 0.0601763686530047

# Row 6 (Neutrophils / Nrg4 / Erbb4 / MuSCs)
$ws.Range("G6").Value = 1.432006666666667
$ws.Range("H6").Value = 4.29602
$ws.Range("I6").Value = 0.1696585030344878
$ws.Range("J6").Value = 0.1846233767260712
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.0108025
$ws.Range("N6").Value = 0.021605
$ws.Range("Q6").Value = 0.01546925201666667
$ws.Range("R6").Value = 0.09281551210000001
$ws.Range("S6").Value = 0.1696585030344878
$ws.Range("T6").Value = 0.1846233767260712

# Row 7 (Resolving-Mac / Nrg4 / Erbb4 / MuSCs)
$ws.Range("G7").Value = 3.034346333333334
$ws.Range("H7").Value = 9.103039000000001
$ws.Range("I7").Value = 0.3594973882348221
$ws.Range("J7").Value = 0.3912071635255698
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.0108025
$ws.Range("N7").Value = 0.021605
$ws.Range("Q7").Value = 0.03277852626583334
$ws.Range("R7").Value = 0.196671157595
$ws.Range("S7").Value = 0.3594973882348221
$ws.Range("T7").Value = 0.3912071635255698
